$wb = $excel.ActiveWorkbook
$wsRef = $wb.Worksheets.Item(1)

# ----------------------------------------------------------------------
# Build a one-time style "palette" on the first worksheet (scratch area,
# row 80) BEFORE touching any real content. Re-using a single already
# finalized source range for every subsequent Copy/PasteSpecial keeps
# the shared style table (styles.xml) clean - mutating alignment/wrap
# properties independently on every sheet can otherwise leave behind
# duplicate/orphan cellXfs entries.
# ----------------------------------------------------------------------
$wsRef.Range("A8").Copy()
$wsRef.Range("A80").PasteSpecial(-4122)          # s6 : time-slot style
$wsRef.Range("B8").Copy()
$wsRef.Range("B80").PasteSpecial(-4122)          # s7 : plain bordered style
$wsRef.Range("D8").Copy()
$wsRef.Range("C80").PasteSpecial(-4122)          # s8 : filled style

$wsRef.Range("B8").Copy()
$wsRef.Range("D80").PasteSpecial(-4122)          # new s9: blank bordered, no alignment
$wsRef.Range("D80").HorizontalAlignment = 1
$wsRef.Range("D80").VerticalAlignment = -4107
$wsRef.Range("D80").WrapText = $false

$timeStyleSrc = $wsRef.Range("A80")
$plainStyleSrc = $wsRef.Range("B80")
$filledStyleSrc = $wsRef.Range("C80")
$blankBorderStyleSrc = $wsRef.Range("D80")

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # ------------------------------------------------------------------
    # Row 8: S1(07:00-09:00)/R101 -> S2(09:00-11:00)/R102
    # The filled "class info" cell moves from D8 to G8, class CL05 -> CL10
    # ------------------------------------------------------------------
    $timeStyleSrc.Copy()
    $ws.Range("A8").PasteSpecial(-4122)
    $ws.Range("A8").Value2 = "S2`n(09:00-11:00)"

    $plainStyleSrc.Copy()
    $ws.Range("B8").PasteSpecial(-4122)
    $ws.Range("B8").Value2 = "R102"

    $plainStyleSrc.Copy()
    $ws.Range("C8:F8").PasteSpecial(-4122)
    $ws.Range("C8:F8").Value2 = ""

    $filledStyleSrc.Copy()
    $ws.Range("G8").PasteSpecial(-4122)
    $ws.Range("G8").Value2 = "Lớp: CL10`nMôn: Tiếng Anh chuyên ngành`nGV: Võ Văn F`n(Lý thuyết)"

    $plainStyleSrc.Copy()
    $ws.Range("H8").PasteSpecial(-4122)
    $ws.Range("H8").Value2 = ""

    $ws.Rows.Item(8).RowHeight = 60

    # ------------------------------------------------------------------
    # Row 9: drop the C2(15:00-17:00) time slot entirely (blank/no time),
    # room becomes R105, and the "class info" moves from G9 to H9, with
    # its text updated to the CL05 / Tieng Anh / Vo Van F entry.
    # ------------------------------------------------------------------
    $blankBorderStyleSrc.Copy()
    $ws.Range("A9").PasteSpecial(-4122)
    $ws.Range("A9").Value2 = ""

    $plainStyleSrc.Copy()
    $ws.Range("B9").PasteSpecial(-4122)
    $ws.Range("B9").Value2 = "R105"

    $plainStyleSrc.Copy()
    $ws.Range("C9:G9").PasteSpecial(-4122)
    $ws.Range("C9:G9").Value2 = ""

    $filledStyleSrc.Copy()
    $ws.Range("H9").PasteSpecial(-4122)
    $ws.Range("H9").Value2 = "Lớp: CL05`nMôn: Tiếng Anh chuyên ngành`nGV: Võ Văn F`n(Lý thuyết)"

    $ws.Rows.Item(9).RowHeight = 15

    # ------------------------------------------------------------------
    # Row 10 (new): C1(13:00-15:00)/R103, carries the former "Ky nang
    # mem / CL05 / Hoang Thi E" entry (this used to be on row 9 / col G).
    # ------------------------------------------------------------------
    $timeStyleSrc.Copy()
    $ws.Range("A10").PasteSpecial(-4122)
    $ws.Range("A10").Value2 = "C1`n(13:00-15:00)"

    $plainStyleSrc.Copy()
    $ws.Range("B10").PasteSpecial(-4122)
    $ws.Range("B10").Value2 = "R103"

    $plainStyleSrc.Copy()
    $ws.Range("C10:E10").PasteSpecial(-4122)
    $ws.Range("C10:E10").Value2 = ""

    $filledStyleSrc.Copy()
    $ws.Range("F10").PasteSpecial(-4122)
    $ws.Range("F10").Value2 = "Lớp: CL05`nMôn: Kỹ năng mềm`nGV: Hoàng Thị E`n(Lý thuyết)"

    $plainStyleSrc.Copy()
    $ws.Range("G10:H10").PasteSpecial(-4122)
    $ws.Range("G10:H10").Value2 = ""

    $ws.Rows.Item(10).RowHeight = 60

    # ------------------------------------------------------------------
    # Row 11 (new): T1(17:30-19:30)/R104, brand-new CL10 / Ky nang mem /
    # Ngo Van I entry.
    # ------------------------------------------------------------------
    $timeStyleSrc.Copy()
    $ws.Range("A11").PasteSpecial(-4122)
    $ws.Range("A11").Value2 = "T1`n(17:30-19:30)"

    $plainStyleSrc.Copy()
    $ws.Range("B11").PasteSpecial(-4122)
    $ws.Range("B11").Value2 = "R104"

    $plainStyleSrc.Copy()
    $ws.Range("C11:F11").PasteSpecial(-4122)
    $ws.Range("C11:F11").Value2 = ""

    $filledStyleSrc.Copy()
    $ws.Range("G11").PasteSpecial(-4122)
    $ws.Range("G11").Value2 = "Lớp: CL10`nMôn: Kỹ năng mềm`nGV: Ngô Văn I`n(Lý thuyết)"

    $plainStyleSrc.Copy()
    $ws.Range("H11").PasteSpecial(-4122)
    $ws.Range("H11").Value2 = ""

    $ws.Rows.Item(11).RowHeight = 60
}

# ----------------------------------------------------------------------
# Clean up the scratch area so the (first) sheet's dimension shrinks
# back down to the real used range.
# ----------------------------------------------------------------------
$wsRef.Range("A80:D80").Clear()
